$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "11222863893"
$ws.Range("B2").Value = "Veenasingh@722"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").VerticalAlignment = -4160
$ws.Range("B2").Locked = $False
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Veenasingh@722", "", "", "Veenasingh@722")
